$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New results data for rows 11-14 (study subjects 8-11), columns B:I
$data = @(
    @(99.2, 97.6, 97.6, 97.540999999999997, 100, 98.4375, 98.461500000000001, 93.75),
    @(93.700800000000001, 91.406300000000002, 98.412700000000001, 98.4375, 100, 92.126000000000005, 99.180300000000003, 99.212599999999995),
    @(94.262299999999996, 89.843800000000002, 89.344300000000004, 89.843800000000002, 98.333299999999994, 98.425200000000004, 94.166700000000006, 97.637799999999999),
    @(95.2, 98.319299999999998, 91.2, 94.117599999999996, 99.206299999999999, 100, 96.031700000000001, 100)
)

$startRow = 11
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 2 + $j  # column B is 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# Update the active selection to B15
$ws.Range("B15").Select()
